$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh crypto market data (price + 1h volume-change column) as scraped on
# Sat Oct 26 18:46:32 UTC 2024. Rows 33/34 also swap rank order (Fetch.AI now
# ranks above PancakeSwap), so those two rows get Coin/Link/Price/Volume all
# rewritten rather than just Price/Volume.
#
# Price column (D) values are plain text (e.g. "67.083.74", "0.0₃0897") that
# must stay text -- force text format before assigning so Excel does not
# coerce them into numbers and lose trailing zeros / the literal formatting.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.083.74"
$ws.Range("E2").Value = "  -0.47%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.481.21"
$ws.Range("E3").Value = "  -0.81%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "584.50"
$ws.Range("E5").Value = "  -0.61%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "170.13"
$ws.Range("E6").Value = "  -0.13%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("E8").Value = "  -1.10%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.480.76"
$ws.Range("E9").Value = "  -0.75%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.136"
$ws.Range("E10").Value = "  +0.89%  "
$ws.Range("E11").Value = "  -0.78%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.93"
$ws.Range("E12").Value = "  -0.29%  "
$ws.Range("E13").Value = "  -2.42%  "
$ws.Range("E14").Value = "  -1.92%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "25.33"
$ws.Range("E15").Value = "  -2.53%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.931.54"
$ws.Range("E16").Value = "  -0.58%  "
$ws.Range("E17").Value = "  -2.03%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.504.83"
$ws.Range("E18").Value = "  +0.53%  "
$ws.Range("E19").Value = "  -6.51%  "
$ws.Range("E20").Value = "  -7.10%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "349.66"
$ws.Range("E21").Value = "  -4.05%  "
$ws.Range("E22").Value = "  -1.44%  "
$ws.Range("E23").Value = "  -0.07%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "68.53"
$ws.Range("E24").Value = "  -4.10%  "
$ws.Range("E25").Value = "  -5.52%  "
$ws.Range("E26").Value = "  -2.58%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.29"
$ws.Range("E27").Value = "  -1.69%  "
$ws.Range("E28").Value = "  -2.34%  "
$ws.Range("E29").Value = "  -0.92%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0897"
$ws.Range("E30").Value = "  -3.93%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "507.88"
$ws.Range("E31").Value = "  -2.00%  "
$ws.Range("E32").Value = "  -6.00%  "
$ws.Range("B33").Value = "Fetch.AI"
$ws.Range("C33").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.23"
$ws.Range("E33").Value = "  -3.48%  "
$ws.Range("B34").Value = "PancakeSwap"
$ws.Range("C34").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.76"
$ws.Range("E34").Value = "  -4.22%  "
$ws.Range("E35").Value = "  -0.11%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "158.72"
$ws.Range("E36").Value = "  +1.63%  "
$ws.Range("E37").Value = "  -8.53%  "
$ws.Range("E38").Value = "  +0.40%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.22"
$ws.Range("E39").Value = "  -4.62%  "
$ws.Range("E40").Value = "  -6.43%  "
$ws.Range("E41").Value = "  -0.33%  "
$ws.Range("E42").Value = "  -4.30%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.326"
$ws.Range("E43").Value = "  -2.01%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.78"
$ws.Range("E44").Value = "  -4.06%  "
$ws.Range("E45").Value = "  -5.13%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "38.73"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "142.39"
$ws.Range("E47").Value = "  -0.91%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.513"
$ws.Range("E48").Value = "  -4.58%  "
$ws.Range("E49").Value = "  -5.47%  "
$ws.Range("E50").Value = "  -6.26%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0728"
$ws.Range("E51").Value = "  -1.49%  "
